# Applies the "Fruta / hortaliza, semanal" update:
# A new weekly price record (row) is inserted at row 11, pushing the
# existing rows 11-41 down to 12-42.
#
# The new row 11 keeps the same market/category/quality/unit/origin
# metadata as the (old) row 11 but carries a new date and new prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 11; this shifts rows 11-41
# down to 12-42 and carries their formatting (incl. the date format on
# column D) along with them.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new record's data.
$ws.Cells.Item(11, 1).Value  = 1
$ws.Cells.Item(11, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(11, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(11, 4).Value  = 44481
$ws.Cells.Item(11, 5).Value  = 15
$ws.Cells.Item(11, 6).Value  = 100112012
$ws.Cells.Item(11, 7).Value  = "Espinaca"
$ws.Cells.Item(11, 8).Value  = "Sin especificar"
$ws.Cells.Item(11, 9).Value  = "Primera"
$ws.Cells.Item(11, 10).Value = 250
$ws.Cells.Item(11, 11).Value = 900
$ws.Cells.Item(11, 12).Value = 1000
$ws.Cells.Item(11, 13).Value = 950
$ws.Cells.Item(11, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(11, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(11, 16).Value = 317
$ws.Cells.Item(11, 17).Value = 3
$ws.Cells.Item(11, 18).Value = "Hortaliza"

# Make sure the date cell keeps the expected date/time number format
# used throughout column D.
$ws.Cells.Item(11, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
